$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AD1:AF1, matching the style of AC1 (bold header)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108
$ws.Range("AD1:AF1").VerticalAlignment = -4160
$ws.Range("AD1:AF1").Borders.LineStyle = 1
$ws.Range("AD1:AF1").Borders.Weight = 2

# Data rows 2-48: team record values
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
